$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 807.6077450938764
$ws.Range("B2").Value = 1156.06823400084
$ws.Range("C2").Value = 2076.469173250837
$ws.Range("D2").Value = 6478.50722779852
$ws.Range("E2").Value = 16215.62762482892
$ws.Range("F2").Value = 51584.20364845017

# Row 3 (A3 unchanged)
$ws.Range("B3").Value = 9.001688094533295
$ws.Range("C3").Value = 8.845335456611505
$ws.Range("D3").Value = 8.609078161932615
$ws.Range("E3").Value = 8.801139489036734
$ws.Range("F3").Value = 10.39958739790717

# Row 5 (A5 unchanged)
$ws.Range("B5").Value = 0.8620551779288285
$ws.Range("C5").Value = 0.6441423430627748
$ws.Range("D5").Value = 0.3682327069172331
$ws.Range("E5").Value = 0.2402838864454218
$ws.Range("F5").Value = 0.1462215113954418

# Row 6 (A6 unchanged)
$ws.Range("B6").Value = 0.00172229409573552
$ws.Range("C6").Value = 0.001287409059188732
$ws.Range("D6").Value = 0.0007352468500899975
$ws.Range("E6").Value = 0.0004788693853570595
$ws.Range("F6").Value = 0.0002902348080217676
